$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.642.15"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.843.67"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'312.46"
$ws.Range("E5").Value = "  -0.89%  "

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4265"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").Value = "'0.3614"
$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("D9").Value = "'0.07288"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").Value = "'0.8699"
$ws.Range("E10").Value = "  -1.95%  "

$ws.Range("D11").Value = "'20.67"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "1.848.64"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").Value = "'6.557"
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").Value = "'5.331"

$ws.Range("D15").Value = "'0.06960"
$ws.Range("E15").Value = "  +1.21%  "

$ws.Range("E16").Value = "  +0.22%  "

$ws.Range("D17").Value = "'79.44"
$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").Value = "'0.000008957"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("E19").Value = "  +0.11%  "

$ws.Range("D20").Value = "'15.29"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").Value = "27.720.33"
$ws.Range("E21").Value = "  +0.34%  "

$ws.Range("D22").Value = "'4.979"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'10.35"
$ws.Range("E23").Value = "  -1.88%  "

$ws.Range("D24").Value = "2.096.49"
$ws.Range("E24").Value = "  +1.59%  "

$ws.Range("D25").Value = "'1.981"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("D26").Value = "'155.12"
$ws.Range("E26").Value = "  -0.26%  "

$ws.Range("D27").Value = "'18.53"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("D28").Value = "'120.53"
$ws.Range("E28").Value = "  -2.09%  "

$ws.Range("D29").Value = "'5.237"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").Value = "'1.870"
$ws.Range("E30").Value = "  +1.51%  "

$ws.Range("D31").Value = "'0.08911"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").Value = "'0.7660"
$ws.Range("E32").Value = "  -1.03%  "

$ws.Range("D33").Value = "'2.967"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("D34").Value = "'4.495"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("E35").Value = "  +2.98%  "

$ws.Range("D36").Value = "'1.001"
$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("D37").Value = "'0.05427"
$ws.Range("E37").Value = "  +1.24%  "

$ws.Range("D38").Value = "'1.103"

$ws.Range("D39").Value = "'0.01926"

$ws.Range("D40").Value = "'2.816"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").Value = "'0.5064"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").Value = "'6.566"
$ws.Range("E43").Value = "  -4.20%  "

$ws.Range("D44").Value = "'8.404"
$ws.Range("E44").Value = "  +1.78%  "

$ws.Range("D45").Value = "'0.06547"
$ws.Range("E45").Value = "  -0.61%  "

$ws.Range("D46").Value = "'106.27"
$ws.Range("E46").Value = "  +1.51%  "

$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "'0.4639"
$ws.Range("E48").Value = "  -1.51%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.001"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").Value = "'1.633"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "'1.752"
$ws.Range("E51").Value = "  -0.64%  "
